# Applies the "Renamed few transcripts. Updated the DataSheet" edit:
#  - Column D (Speaker) values are shortened:
#       "CECILIO DIMAS" -> "T"
#       "STUDENT"       -> "S"
#  - Column F (Teacher Tag) value "3 - Getting Students to Relate" is
#    shortened to "3 - Getting SS to Relate" wherever it occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $speakerCell = $ws.Cells.Item($r, 4)
    $speaker = $speakerCell.Text

    if ($speaker -eq "CECILIO DIMAS") {
        $speakerCell.Value = "T"
    } elseif ($speaker -eq "STUDENT") {
        $speakerCell.Value = "S"
    }

    $tagCell = $ws.Cells.Item($r, 6)
    $tag = $tagCell.Text

    if ($tag -eq "3 - Getting Students to Relate") {
        $tagCell.Value = "3 - Getting SS to Relate"
    }
}
